# Actualización automática de grupos experimentales
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Toggle / set "Grupo_Experimental" (column B) values for the affected rows
$ws.Range("B2").Value = "Sin SmartScore"
$ws.Range("B3").Value = "Con SmartScore"
$ws.Range("B4").Value = "Con SmartScore"
$ws.Range("B7").Value = "Con SmartScore"
$ws.Range("B8").Value = "Con SmartScore"
$ws.Range("B10").Value = "Sin SmartScore"
$ws.Range("B11").Value = "Sin SmartScore"
$ws.Range("B12").Value = "Con SmartScore"
$ws.Range("B13").Value = "Sin SmartScore"
$ws.Range("B16").Value = "Sin SmartScore"

# Convert the SmartScore numeric columns in row 16 from text to real numbers
$ws.Range("I16").Value = 0.549
$ws.Range("L16").Value = 0.526
$ws.Range("O16").Value = 0.491
$ws.Range("R16").Value = 0.612
$ws.Range("U16").Value = 0.52
$ws.Range("X16").Value = 0.518
$ws.Range("AA16").Value = 0.725
$ws.Range("AD16").Value = 0.605
$ws.Range("AG16").Value = 0.571
